$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) -----------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value  = 566
$wsExhibit.Range("F11").Value = 2799
$wsExhibit.Range("F14").Value = 1101
$wsExhibit.Range("F18").Value = 1604
$wsExhibit.Range("F25").Value = 1467
$wsExhibit.Range("F26").Value = 1455
$wsExhibit.Range("F28").Value = 270
$wsExhibit.Range("F39").Value = 2274
$wsExhibit.Range("F42").Value = 2767

# --- Sheet "演出" (Performances) -----------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F12").Value = 364
$wsShow.Range("F13").Value = 0

# --- Sheet "本地生活" (Local Life) ---------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F13").Value = 1347

# --- Sheet "全部类型" (All Types) ----------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

# The "start date" column (B) holds plain text like "2024-06-01" in this
# workbook, not real dates. Force text format first so Excel doesn't
# auto-convert the assigned string into a date serial number.
# (B8's date text is unchanged by this edit, so it is left untouched.)
$wsAll.Range("B5:B7").NumberFormat = "@"

# Row 5 now holds what used to be the row-6 event (NIJISANJI EN)
$wsAll.Range("B5").Value = "2024-06-01"
$wsAll.Range("C5").Value = "上海·NIJISANJI EN 官方授权主题店"
$wsAll.Range("D5").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$wsAll.Range("E5").Value = "2024.06.01 00:00-07.15 23:59"
$wsAll.Range("F5").Value = 663
$wsAll.Range("G5").Value = 30
$wsAll.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=86310"
$wsAll.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202405/MhBVkfZ51716778059321.jpeg"

# Row 6 now holds what used to be the row-7 event (全职高手 x HAPPY ZOO)
$wsAll.Range("B6").Value = "2024-06-07"
$wsAll.Range("C6").Value = "上海·全职高手×HAPPY ZOO 全职高手十周年咖啡厅"
$wsAll.Range("D6").Value = "南京东路340号百联zx创趣场四楼05号 HAPPY ZOO"
$wsAll.Range("E6").Value = "2024.06.07 00:00-08.04 23:59"
$wsAll.Range("F6").Value = 931
$wsAll.Range("G6").Value = 10
$wsAll.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86871"
$wsAll.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/KLJmCEkC1717568198482.png"

# Row 7 now holds what used to be the row-8 event (怪兽8号 meets niko and …)
$wsAll.Range("B7").Value = "2024-06-08"
$wsAll.Range("C7").Value = "上海· 怪兽8号 meets niko and … 集章之旅    "
$wsAll.Range("D7").Value = "吴江路169号1层E127,E128 niko and ... (上海四季坊店)"
$wsAll.Range("E7").Value = "2024.06.08 10:00-07.21 22:00"
$wsAll.Range("F7").Value = 541
$wsAll.Range("G7").Value = 48
$wsAll.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=85758"
$wsAll.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202405/xw8aUE5u1715846379865.jpeg"

# Row 8 is a brand new event (无穹-中国 航天沉浸艺术展).
# B8 already reads "2024-06-08" and is unchanged by this edit, so it is
# intentionally left alone (re-assigning it would risk Excel re-parsing
# the text as a date serial).
$wsAll.Range("C8").Value = "上海·无穹-中国 航天沉浸艺术展"
$wsAll.Range("D8").Value = "上海浦东新区樱花路869号3F 上海喜玛拉雅美术馆"
$wsAll.Range("E8").Value = "2024.06.08 10:00-10.07 20:00"
$wsAll.Range("F8").Value = 101
$wsAll.Range("G8").Value = 78
$wsAll.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=86957"
$wsAll.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202406/Bus3lAnI1717558639134.jpeg"

# Remaining simple numeric bumps on "全部类型"
$wsAll.Range("F10").Value = 1347
$wsAll.Range("F16").Value = 566
$wsAll.Range("F18").Value = 2799
$wsAll.Range("F22").Value = 1101
$wsAll.Range("F25").Value = 1604
$wsAll.Range("F27").Value = 364
$wsAll.Range("F31").Value = 1467
$wsAll.Range("F32").Value = 1455
$wsAll.Range("F43").Value = 2274
$wsAll.Range("F46").Value = 2767
